$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-column template row used to copy the existing cell style (number
# format) onto freshly added cells, so the new cells reuse the existing
# style indices instead of the engine allocating brand-new ones.
$fmtSrcRow = 470
$fmtSrcRowN = 467   # row 470 has no populated "N" cell to copy from

$newRows = @(
    @{ row=471; A=45193.82327327546; B="dlruddk9@naver.com"; C="사회복지학과"; D=20212342; E="이경아"; F="‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."; G=0.5; H="5:5"; I="20분의 1"; J="44만호, 153만명"; K="전라"; L="Red"; M="반대한다." },
    @{ row=472; A=45193.824795555556; B="ayden0429@gmail.com"; C="의예과"; D=20226145; E="이성연"; F="‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."; G=0.1; H="6:4"; I="20분의 1"; J="20만호, 69만명"; K="충청"; L="Red"; M="근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다." },
    @{ row=473; A=45193.830308414355; B="lcbat4@gmail.com"; C="글로벌학부"; D=20236429; E="홍서경"; F="과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."; G=0.7; H="6:4"; I="20분의 1"; J="20만호, 69만명"; K="경상"; L="Black"; N="모름/무응답" },
    @{ row=474; A=45193.832151145834; B="hkmcosmos1@gmail.com"; C="글로벌비즈니스"; D=20226429; E="한기민"; F="‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."; G=0.1; H="3:7"; I="15분의 1"; J="20만호, 69만명"; K="충청"; L="Black"; N="모름/무응답" },
    @{ row=475; A=45193.837822326386; B="kimbitna7890@naver.com"; C="광고홍보학과"; D=20222609; E="김빛나"; F="과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."; G=0.3; H="3:7"; I="15분의 1"; J="15만호,  32만명"; K="평안"; L="Black"; N="노동자가 과도한 연장근로를 받을 수 있어 반대한다." },
    @{ row=476; A=45193.83826665509; B="hyelinj27@gmail.com"; C="인공지능융합학부"; D=20236781; E="진혜린"; F="등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."; G=0.1; H="6:4"; I="20분의 1"; J="20만호, 69만명"; K="충청"; L="Black"; N="찬성한다." },
    @{ row=477; A=45193.84510828704; B="sillysunny@naver.com"; C="인문학부"; D=20231037; E="박세현"; F="‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."; G=0.7; H="6:4"; I="20분의 1"; J="20만호, 69만명"; K="전라"; L="Red"; M="모름/무응답" },
    @{ row=478; A=45193.852299907405; B="ljh2017@naver.com"; C="철학과"; D=20181079; E="이정효"; F="등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."; G=0.3; H="3:7"; I="10분의 1"; J="44만호, 153만명"; K="전라"; L="Black"; N="모름/무응답" },
    @{ row=479; A=45193.85701881944; B="jisung5549@naver.com"; C="경영학과"; D=20222970; E="송지성"; F="실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."; G=0.7; H="4:6"; I="20분의 1"; J="44만호, 153만명"; K="경상"; L="Red"; M="근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다." },
    @{ row=480; A=45193.862970636575; B="jangho5636@gmail.com"; C="러시아학과"; D=20161723; E="이장호"; F="등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."; G=0.9; H="7:3"; I="30분의 1"; J="20만호, 69만명"; K="평안"; L="Red"; M="반대한다." },
    @{ row=481; A=45193.8716166088; B="0227jsh@naver.com"; C="식품영양학과"; D=20233843; E="장서희"; F="과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."; G=0.3; H="6:4"; I="20분의 1"; J="15만호,  32만명"; K="평안"; L="Red"; M="근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다." },
    @{ row=482; A=45193.8724659375; B="rhksan324@naver.com"; C="금융재무학과"; D=20203001; E="이관무"; F="등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."; G=0.1; H="6:4"; I="20분의 1"; J="20만호, 69만명"; K="충청"; L="Red"; M="근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다." },
    @{ row=483; A=45193.8728371412; B="jihye199530@gmail.com"; C="간호학과"; D=20236261; E="엄지혜"; F="‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."; G=0.7; H="5:5"; I="20분의 1"; J="15만호,  32만명"; K="평안"; L="Black"; N="노동자가 과도한 연장근로를 받을 수 있어 반대한다." },
    @{ row=484; A=45193.87713403935; B="jongbaep17s@gmail.com"; C="글로벌비즈니스"; D=20226410; E="박종배"; F="‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."; G=0.1; H="6:4"; I="20분의 1"; J="20만호, 69만명"; K="충청"; L="Red"; M="모름/무응답" },
    @{ row=485; A=45193.9036634838; B="rhdskrud123@naver.com"; C="인문학부"; D=20231002; E="공나경"; F="실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."; G=0.7; H="7:3"; I="20분의 1"; J="20만호, 69만명"; K="전라"; L="Black"; N="노동자가 과도한 연장근로를 받을 수 있어 반대한다." },
    @{ row=486; A=45193.90699703703; B="shdbsgh0305@naver.com"; C="러시아학과"; D=20231710; E="노윤호"; F="‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."; G=0.5; H="5:5"; I="20분의 1"; J="20만호, 69만명"; K="전라"; L="Black"; N="모름/무응답" },
    @{ row=487; A=45193.90846236111; B="digiphk12@naver.com"; C="광고홍보학과"; D=20162617; E="박현규"; F="등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."; G=0.7; H="6:4"; I="15분의 1"; J="20만호, 69만명"; K="평안"; L="Red"; M="근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다." },
    @{ row=488; A=45193.909824409726; B="mnsghn314@naver.com"; C="소프트웨어학과"; D=20235159; E="문승현"; F="등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."; G=0.1; H="6:4"; I="10분의 1"; J="20만호, 69만명"; K="충청"; L="Red"; M="근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다." },
    @{ row=489; A=45193.91182700232; B="amielee1997@naver.com"; C="간호학과"; D=20217159; E="이지수"; F="과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."; G=0.1; H="7:3"; I="10분의 1"; J="15만호,  32만명"; K="경기"; L="Black"; N="모름/무응답" },
    @{ row=490; A=45193.91356240741; B="sinfkks@gmail.com"; C="반도체 디스플레이스쿨"; D=20233304; E="김경진"; F="‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."; G=0.9; H="4:6"; I="10분의 1"; J="44만호, 153만명"; K="경기"; L="Black"; N="찬성한다." }
)

foreach ($r in $newRows) {
    $rn = $r["row"]
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N")) {
        if ($r.ContainsKey($col)) {
            $srcRow = $fmtSrcRow
            if ($col -eq "N") { $srcRow = $fmtSrcRowN }
            $ws.Range("${col}${srcRow}").Copy()
            $ws.Range("${col}${rn}").PasteSpecial(-4122)
            $ws.Range("${col}${rn}").Value = $r[$col]
        }
    }
}

$excel.CutCopyMode = 0